$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -3
    3  = -1
    5  = 4
    7  = 3
    8  = -6
    9  = 3
    10 = -3
    11 = -2
    12 = -6
    13 = 4
    15 = -1
    16 = -3
    17 = -1
    19 = -2
    20 = -4
    21 = -5
    22 = -1
    23 = -1
    24 = -5
    25 = 1
    27 = -3
    28 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
